# Update shadow model values on Sheet1 (column D) according to
# D1/L1 = D2/L2 --> D2 = ... recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = 2926.799320310989
    "D3"  = 2926.799320310986
    "D4"  = 2926.799320310986
    "D6"  = 8362.283772317136
    "D7"  = 8362.283772317136
    "D9"  = 9168.262111917831
    "D10" = 9168.262111917824
    "D11" = 133.3333333333333
    "D12" = 7568.262111917824
    "D13" = 133.3333333333333
    "D17" = 106666.6666666667
    "D19" = 121802.6639999999
    "D20" = 121802.6639999999
    "D21" = 1600
    "D24" = 143294.1379324235
    "D25" = 143294.1379324234
    "D28" = 2865.882758648496
    "D29" = 2865.882758648496
    "D30" = 143294.1379324234
    "D35" = 20217.59999999971
    "D36" = 20217.59999999971
    "D38" = -5489.676688001595
    "D39" = -5489.676688001589
    "D41" = 5489.676688001589
    "D42" = 109793.5337600332
    "D43" = 109793.5337600332
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
